$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value2 = 909526.9399999999
$ws.Range("I28").Value2 = 2000112.8
$ws.Range("J28").Value2 = 705.3333
$ws.Range("K28").Value2 = 2000112.8
$ws.Range("L28").Value2 = 705.3333
$ws.Range("M28").Value2 = -1999627.8
$ws.Range("N28").Value2 = -1675.3333
$ws.Range("H43").Value2 = 1500
$ws.Range("J43").Value2 = 0
$ws.Range("L43").Value2 = 0
$ws.Range("N43").ClearContents()
$ws.Range("H55").Value2 = 76.666664
$ws.Range("I55").Value2 = 99.818184
$ws.Range("K55").Value2 = 99.818184
$ws.Range("M55").Value2 = 114.181816
$ws.Range("H82").Value2 = 14286332
$ws.Range("I82").Value2 = 16667306
$ws.Range("K82").Value2 = 50001918
$ws.Range("M82").Value2 = -50001512
$ws.Range("H85").Value2 = 14286332
$ws.Range("I85").Value2 = 16667306
$ws.Range("K85").Value2 = 50001918
$ws.Range("M85").Value2 = -50000514
$ws.Range("H124").Value2 = 99994
$ws.Range("J124").Value2 = 99994
$ws.Range("L124").Value2 = 99994
$ws.Range("N124").Value2 = -109814
$ws.Range("H130").Value2 = 87984.5
$ws.Range("J130").Value2 = 87984.5
$ws.Range("L130").Value2 = 87984.5
$ws.Range("N130").Value2 = -98024.5
$ws.Range("H135").Value2 = 3897.9
$ws.Range("J135").Value2 = 1499
$ws.Range("L135").Value2 = 13491
$ws.Range("N135").Value2 = -18561

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 1672.4512
$ws.Range("I32").Value2 = 1172.2894
$ws.Range("J32").Value2 = 8007.8335
$ws.Range("K32").Value2 = 1172.2894
$ws.Range("L32").Value2 = 8007.8335
$ws.Range("M32").Value2 = -885.2893999999999
$ws.Range("N32").Value2 = -8581.833500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value2 = 702.17645
$ws.Range("J80").Value2 = 682.7143
$ws.Range("L80").Value2 = 682.7143
$ws.Range("N80").Value2 = -2678.7143
$ws.Range("H83").Value2 = 702.17645
$ws.Range("J83").Value2 = 682.7143
$ws.Range("L83").Value2 = 3413.5715
$ws.Range("N83").Value2 = -13397.5715
$ws.Range("H122").Value2 = 99890
$ws.Range("J122").Value2 = 99890
$ws.Range("L122").Value2 = 99890
$ws.Range("N122").Value2 = -109690

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value2 = 5400
$ws.Range("J2").Value2 = 5400
$ws.Range("L2").Value2 = 5400
$ws.Range("N2").Value2 = -5626
$ws.Range("H5").Value2 = 1377.7693
$ws.Range("J5").Value2 = 1626.4546
$ws.Range("L5").Value2 = 1626.4546
$ws.Range("N5").Value2 = -1850.4546
$ws.Range("H10").Value2 = 1712.25
$ws.Range("I10").Value2 = 613.6667
$ws.Range("J10").Value2 = 5008
$ws.Range("K10").Value2 = 613.6667
$ws.Range("L10").Value2 = 5008
$ws.Range("M10").Value2 = -474.6667
$ws.Range("N10").Value2 = -5286
$ws.Range("H11").Value2 = 2298.6
$ws.Range("I11").Value2 = 498
$ws.Range("J11").Value2 = 3499
$ws.Range("K11").Value2 = 498
$ws.Range("L11").Value2 = 3499
$ws.Range("M11").Value2 = -358
$ws.Range("N11").Value2 = -3779
$ws.Range("H13").Value2 = 7497.5
$ws.Range("I13").Value2 = 3496
$ws.Range("K13").Value2 = 3496
$ws.Range("M13").Value2 = -3357
$ws.Range("H14").Value2 = 27998.5
$ws.Range("J14").Value2 = 27998.5
$ws.Range("L14").Value2 = 27998.5
$ws.Range("N14").Value2 = -28338.5
$ws.Range("H23").Value2 = 0
$ws.Range("I23").Value2 = 0
$ws.Range("K23").Value2 = 0
$ws.Range("M23").ClearContents()
$ws.Range("H26").Value2 = 28624.75
$ws.Range("I26").Value2 = 5250
$ws.Range("J26").Value2 = 51999.5
$ws.Range("K26").Value2 = 5250
$ws.Range("L26").Value2 = 51999.5
$ws.Range("M26").Value2 = -4963
$ws.Range("N26").Value2 = -52573.5
$ws.Range("H27").Value2 = 0
$ws.Range("I27").Value2 = 0
$ws.Range("K27").Value2 = 0
$ws.Range("M27").ClearContents()
$ws.Range("H31").Value2 = 4079.5952
$ws.Range("I31").Value2 = 1066.5385
$ws.Range("K31").Value2 = 1066.5385
$ws.Range("M31").Value2 = -771.5385000000001
$ws.Range("H32").Value2 = 1925.8
$ws.Range("I32").Value2 = 1925.8
$ws.Range("K32").Value2 = 1925.8
$ws.Range("M32").Value2 = -1609.8
$ws.Range("H33").Value2 = 2022
$ws.Range("I33").Value2 = 2022
$ws.Range("K33").Value2 = 2022
$ws.Range("M33").Value2 = -1643
$ws.Range("H34").Value2 = 4079.5952
$ws.Range("I34").Value2 = 1066.5385
$ws.Range("K34").Value2 = 1066.5385
$ws.Range("M34").Value2 = -864.5385000000001
$ws.Range("H35").Value2 = 988.3077
$ws.Range("I35").Value2 = 250
$ws.Range("J35").Value2 = 2649.5
$ws.Range("K35").Value2 = 250
$ws.Range("L35").Value2 = 2649.5
$ws.Range("M35").Value2 = 44
$ws.Range("N35").Value2 = -3237.5
$ws.Range("H39").Value2 = 14644.4
$ws.Range("I39").Value2 = 14111
$ws.Range("K39").Value2 = 14111
$ws.Range("M39").Value2 = -13720
$ws.Range("H44").Value2 = 43354.668
$ws.Range("I44").Value2 = 45032
$ws.Range("J44").Value2 = 40000
$ws.Range("K44").Value2 = 45032
$ws.Range("L44").Value2 = 40000
$ws.Range("M44").Value2 = -44590
$ws.Range("N44").Value2 = -40884
$ws.Range("H48").Value2 = 0
$ws.Range("J48").Value2 = 0
$ws.Range("L48").Value2 = 0
$ws.Range("N48").ClearContents()
$ws.Range("H49").Value2 = 14644.4
$ws.Range("I49").Value2 = 14111
$ws.Range("K49").Value2 = 14111
$ws.Range("M49").Value2 = -13929
$ws.Range("H59").Value2 = 49500
$ws.Range("J59").Value2 = 49500
$ws.Range("L59").Value2 = 49500
$ws.Range("N59").Value2 = -51790
$ws.Range("H105").Value2 = 1470.0646
$ws.Range("I105").Value2 = 1524.0625
$ws.Range("J105").Value2 = 1412.4667
$ws.Range("K105").Value2 = 1524.0625
$ws.Range("L105").Value2 = 1412.4667
$ws.Range("M105").Value2 = 222.9375
$ws.Range("N105").Value2 = -4906.4667
$ws.Range("H107").Value2 = 1016.5333
$ws.Range("I107").Value2 = 951.4
$ws.Range("J107").Value2 = 1146.8
$ws.Range("K107").Value2 = 951.4
$ws.Range("L107").Value2 = 1146.8
$ws.Range("M107").Value2 = 968.6
$ws.Range("N107").Value2 = -4986.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 68751010
$ws.Range("I4").Value2 = 100000610
$ws.Range("J4").Value2 = 1898.2
$ws.Range("K4").Value2 = 300001830
$ws.Range("L4").Value2 = 5694.6
$ws.Range("M4").Value2 = -300001718
$ws.Range("N4").Value2 = -5918.6
$ws.Range("H74").Value2 = 20000
$ws.Range("J74").Value2 = 20000
$ws.Range("L74").Value2 = 60000
$ws.Range("N74").Value2 = -62122
$ws.Range("H77").Value2 = 20000
$ws.Range("J77").Value2 = 20000
$ws.Range("L77").Value2 = 180000
$ws.Range("N77").Value2 = -190608
$ws.Range("H131").Value2 = 2741435.5
$ws.Range("I131").Value2 = 6667816.5
$ws.Range("J131").Value2 = 1725992.1
$ws.Range("K131").Value2 = 20003449.5
$ws.Range("L131").Value2 = 5177976.300000001
$ws.Range("M131").Value2 = -19998409.5
$ws.Range("N131").Value2 = -5188056.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 50989
$ws.Range("I70").Value2 = 70494.875
$ws.Range("J70").Value2 = 19779.6
$ws.Range("K70").Value2 = 70494.875
$ws.Range("L70").Value2 = 19779.6
$ws.Range("M70").Value2 = -70224.875
$ws.Range("N70").Value2 = -20319.6
$ws.Range("H73").Value2 = 50989
$ws.Range("I73").Value2 = 70494.875
$ws.Range("J73").Value2 = 19779.6
$ws.Range("K73").Value2 = 70494.875
$ws.Range("L73").Value2 = 19779.6
$ws.Range("M73").Value2 = -69558.875
$ws.Range("N73").Value2 = -21651.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 31255318
$ws.Range("I7").Value2 = 5558.3
$ws.Range("K7").Value2 = 5558.3
$ws.Range("M7").Value2 = -5446.3
$ws.Range("H126").Value2 = 31255318
$ws.Range("I126").Value2 = 5558.3
$ws.Range("K126").Value2 = 16674.9
$ws.Range("M126").Value2 = -14204.9
